# Insert a new weekly record at row 3 (pushing all existing data rows down
# by one) for "Vega Monumental Concepción - Pepino ensalada".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 3..109 down to 4..110, making room for the new record.
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the latest weekly price data.
$ws.Cells.Item(3, 1).Value = 11
$ws.Cells.Item(3, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value = "Bíobío"
$ws.Cells.Item(3, 4).Value = 44631
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 100112043
$ws.Cells.Item(3, 7).Value = "Pepino ensalada"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 220
$ws.Cells.Item(3, 11).Value = 17000
$ws.Cells.Item(3, 12).Value = 18000
$ws.Cells.Item(3, 13).Value = 17545
$ws.Cells.Item(3, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(3, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(3, 16).Value = 292
$ws.Cells.Item(3, 17).Value = 60
$ws.Cells.Item(3, 18).Value = "Hortaliza"
